# Auto-generated Excel COM-interop script
# Applies updated Leve profit-calculation values across all 8 class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 481.625
$ws.Range("I33").Value = 303.45
$ws.Range("J33").Value = 1372.5
$ws.Range("K33").Value = 303.45
$ws.Range("L33").Value = 1372.5
$ws.Range("M33").Value = -74.44999999999999
$ws.Range("N33").Value = -1830.5
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
# Row 95
$ws.Range("H95").Value = 37925
$ws.Range("J95").Value = 37925
$ws.Range("L95").Value = 37925
$ws.Range("N95").Value = -43417
# Row 105
$ws.Range("H105").Value = 49331.5
$ws.Range("J105").Value = 49331.5
$ws.Range("L105").Value = 49331.5
$ws.Range("N105").Value = -56319.5
# Row 108
$ws.Range("H108").Value = 45659
$ws.Range("J108").Value = 45659
$ws.Range("L108").Value = 45659
$ws.Range("N108").Value = -53339
# Row 109
$ws.Range("H109").Value = 39982.668
$ws.Range("J109").Value = 39982.668
$ws.Range("L109").Value = 39982.668
$ws.Range("N109").Value = -42756.668
# Row 130
$ws.Range("H130").Value = 41728
$ws.Range("J130").Value = 41728
$ws.Range("L130").Value = 41728
$ws.Range("N130").Value = -51768
# Row 138
$ws.Range("H138").Value = 1676.1086
$ws.Range("I138").Value = 590.6667
$ws.Range("J138").Value = 2860.2273
$ws.Range("K138").Value = 1772.0001
$ws.Range("L138").Value = 8580.6819
$ws.Range("M138").Value = 3367.9999
$ws.Range("N138").Value = -18860.6819
# Row 141
$ws.Range("H141").Value = 1240
$ws.Range("I141").Value = 695.9375
$ws.Range("J141").Value = 6463
$ws.Range("K141").Value = 2087.8125
$ws.Range("L141").Value = 19389
$ws.Range("M141").Value = 3092.1875
$ws.Range("N141").Value = -29749

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12324.23
$ws.Range("I32").Value = 13489.383
$ws.Range("J32").Value = 9281.888999999999
$ws.Range("K32").Value = 13489.383
$ws.Range("L32").Value = 9281.888999999999
$ws.Range("M32").Value = -13202.383
$ws.Range("N32").Value = -9855.888999999999
# Row 80
$ws.Range("H80").Value = 48715.555
$ws.Range("J80").Value = 48715.555
$ws.Range("L80").Value = 48715.555
$ws.Range("N80").Value = -50711.555
# Row 83
$ws.Range("H83").Value = 48715.555
$ws.Range("J83").Value = 48715.555
$ws.Range("L83").Value = 146146.665
$ws.Range("N83").Value = -156130.665
# Row 101
$ws.Range("H101").Value = 48141.5
$ws.Range("J101").Value = 48141.5
$ws.Range("L101").Value = 48141.5
$ws.Range("N101").Value = -54631.5
# Row 109
$ws.Range("H109").Value = 46559.832
$ws.Range("J109").Value = 46559.832
$ws.Range("L109").Value = 46559.832
$ws.Range("N109").Value = -49333.832
# Row 117
$ws.Range("H117").Value = 46995.5
$ws.Range("J117").Value = 46995.5
$ws.Range("L117").Value = 46995.5
$ws.Range("N117").Value = -56173.5
# Row 118
$ws.Range("H118").Value = 49172
$ws.Range("J118").Value = 49172
$ws.Range("L118").Value = 49172
$ws.Range("N118").Value = -52486

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2517.037
$ws.Range("I20").Value = 1994.625
$ws.Range("J20").Value = 3276.9092
$ws.Range("K20").Value = 1994.625
$ws.Range("L20").Value = 3276.9092
$ws.Range("M20").Value = -1747.625
$ws.Range("N20").Value = -3770.9092

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2818.93
$ws.Range("I31").Value = 1275.2759
$ws.Range("J31").Value = 3449.4365
$ws.Range("K31").Value = 1275.2759
$ws.Range("L31").Value = 3449.4365
$ws.Range("M31").Value = -980.2759000000001
$ws.Range("N31").Value = -4039.4365
# Row 34
$ws.Range("H34").Value = 2818.93
$ws.Range("I34").Value = 1275.2759
$ws.Range("J34").Value = 3449.4365
$ws.Range("K34").Value = 1275.2759
$ws.Range("L34").Value = 3449.4365
$ws.Range("M34").Value = -1073.2759
$ws.Range("N34").Value = -3853.4365
# Row 43
$ws.Range("H43").Value = 46824.5
$ws.Range("J43").Value = 46824.5
$ws.Range("L43").Value = 46824.5
$ws.Range("N43").Value = -47192.5
# Row 101
$ws.Range("H101").Value = 46824.5
$ws.Range("J101").Value = 46824.5
$ws.Range("L101").Value = 46824.5
$ws.Range("N101").Value = -53314.5
# Row 107
$ws.Range("H107").Value = 317.97562
$ws.Range("I107").Value = 248.69444
$ws.Range("J107").Value = 816.8
$ws.Range("K107").Value = 248.69444
$ws.Range("L107").Value = 816.8
$ws.Range("M107").Value = 1671.30556
$ws.Range("N107").Value = -4656.8
# Row 132
$ws.Range("H132").Value = 85308.47
$ws.Range("I132").Value = 2887.25
$ws.Range("J132").Value = 283119.4
$ws.Range("K132").Value = 8661.75
$ws.Range("L132").Value = 849358.2000000001
$ws.Range("M132").Value = -6131.75
$ws.Range("N132").Value = -854418.2000000001

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 139
$ws.Range("H139").Value = 99474.91
$ws.Range("I139").Value = 195948.88
$ws.Range("J139").Value = 3000.9375
$ws.Range("K139").Value = 587846.64
$ws.Range("L139").Value = 9002.8125
$ws.Range("M139").Value = -582706.64
$ws.Range("N139").Value = -19282.8125
# Row 140
$ws.Range("H140").Value = 40906.92
$ws.Range("I140").Value = 54235.844
$ws.Range("K140").Value = 162707.532
$ws.Range("M140").Value = -157527.532

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 11200
$ws.Range("I70").Value = 15666.667
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 15666.667
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -15396.667
$ws.Range("N70").Value = -5040
# Row 73
$ws.Range("H73").Value = 11200
$ws.Range("I73").Value = 15666.667
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 15666.667
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -14730.667
$ws.Range("N73").Value = -6372
# Row 101
$ws.Range("H101").Value = 45653
$ws.Range("J101").Value = 45653
$ws.Range("L101").Value = 45653
$ws.Range("N101").Value = -52143

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 105
$ws.Range("H105").Value = 47303.5
$ws.Range("J105").Value = 47303.5
$ws.Range("L105").Value = 47303.5
$ws.Range("N105").Value = -54291.5

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 34301
$ws.Range("J103").Value = 34301
$ws.Range("L103").Value = 34301
$ws.Range("N103").Value = -36645
# Row 104
$ws.Range("H104").Value = 48361
$ws.Range("J104").Value = 48361
$ws.Range("L104").Value = 48361
$ws.Range("N104").Value = -55349
# Row 119
$ws.Range("H119").Value = 48678
$ws.Range("J119").Value = 48678
$ws.Range("L119").Value = 48678
$ws.Range("N119").Value = -58354
# Row 132
$ws.Range("H132").Value = 2279.8647
$ws.Range("I132").Value = 2113.5186
$ws.Range("J132").Value = 2729
$ws.Range("K132").Value = 6340.5558
$ws.Range("L132").Value = 8187
$ws.Range("M132").Value = -3810.5558
$ws.Range("N132").Value = -13247

# Row 52 on ALC: LeveProfitNQ (M52) is no longer populated for this entry
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("M52").ClearContents()
